# Fed: Update for 2025.03 release
# Applies the resume-content edits described by the target diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Profile / summary blurb ---
Replace-Text "Versatile polyglot programmer proficient in python, SQL (databases), R, etc. Skilled in linux, OOP, and data science. Able to dissect complex problems and refactor into microservice solutions." `
             "Versatile programmer proficient in python, SQL (databases), R, etc. Skilled in Linux, OOP, data science. Teaching experience."

# --- Federal Reserve Bank of Minneapolis: Sr. Data Scientist bullet points ---
Replace-Text "Contributed to infrastructure, econometric models, and dashboards." `
             "Contributed to econometric models, infrastructure, dashboards, adhoc analyses for consumption by Board of Governors (Jerome Powell)."

Replace-Text "Bridged the gap between frontend and database support team, and Production econometric modeling team to ensure smooth integration. Developed database version diff tool to rapidly identify and resolve issues." `
             "Served as Production econometric modeling team’s technical liaison to frontend, database and sysadmin teams, ensuring smooth integration."

Replace-Text "Proactively introduced unit testing, autotesting, and autodocs to python repos. Co-championed long-term plan to streamline code releases." `
             "Technical lead in migrating Production codebase across OS, language versions, environments, while distributing and popularizing DIY automation tools to support fellow quants and economists."

Replace-Text "Independently worked to improve model reproducibility via automation of security, conda, ssh, git, variable OS modules. Worked on OS migration, chmod enablement and GitLab runners." `
             "Proactively introduced unit testing, autotesting, and autodocs to python repos."

Replace-Text "Co-prototyped a fullstack data management system using Flask and SQLite using dynamic SQL queries, HTML forms, endpoints. Currently productionalizing." `
             "Co-prototyped a fullstack data management system using Flask and SQLite using dynamic SQL queries, HTML forms, endpoints."

# --- Medica bullet: reword "reconcile databases" -> "database reconciliation" ---
Replace-Text "Discovered a technique to automate reconcile databases, requiring" `
             "Discovered a technique to automate database reconciliation, requiring"

# --- College/high-school jobs: "business" -> "art studio" ---
Replace-Text "College: Owner of brick-and-mortar business providing website design, search engine optimization and graphic design services." `
             "College: Owner of brick-and-mortar art studio providing website design, search engine optimization and graphic design services."

# --- Skills section styles: re-touch Bold so rPr serializes in canonical (b, bCs) order ---
Replace-Text "Minneapolis, MN -" "Minneapolis, MN -"

# === Style-sheet touch-ups mirrored from the diff ===

# Abstract paragraph style: space-before 300 -> 100 twips (15pt -> 5pt)
$abstractStyle = $d.Styles("Abstract")
$abstractStyle.ParagraphFormat.SpaceBefore = 5

# New "Abstract Title" style
$titleStyle = $d.Styles.Add("AbstractTitle", 1)
$titleStyle.NameLocal = "Abstract Title"
$titleStyle.BaseStyle = $d.Styles("Normal")
$titleStyle.NextParagraphStyle = $d.Styles("Abstract")
$titleStyle.QuickStyle = $true
$titleStyle.ParagraphFormat.KeepWithNext = $true
$titleStyle.ParagraphFormat.KeepTogether = $true
$titleStyle.ParagraphFormat.SpaceBefore = 15
$titleStyle.ParagraphFormat.SpaceAfter = 0
$titleStyle.ParagraphFormat.Alignment = 1
$titleStyle.Font.Bold = $true
$titleStyle.Font.Color = 9067060
$titleStyle.Font.Size = 10
$titleStyle.Font.SizeBi = 10

# New "Footnote Block Text" style
$fnBlockStyle = $d.Styles.Add("FootnoteBlockText", 1)
$fnBlockStyle.NameLocal = "Footnote Block Text"
$fnBlockStyle.BaseStyle = $d.Styles("FootnoteText")
$fnBlockStyle.NextParagraphStyle = $d.Styles("FootnoteText")
$fnBlockStyle.Priority = 9
$fnBlockStyle.UnhideWhenUsed = $true
$fnBlockStyle.QuickStyle = $true
$fnBlockStyle.ParagraphFormat.SpaceBefore = 5
$fnBlockStyle.ParagraphFormat.SpaceAfter = 5
$fnBlockStyle.ParagraphFormat.FirstLineIndent = 0
$fnBlockStyle.ParagraphFormat.LeftIndent = 24
$fnBlockStyle.ParagraphFormat.RightIndent = 24

# Pandoc-highlight character styles: re-assert formatting so element order normalizes,
# and fill in the two previously-blank token styles.
$importTok = $d.Styles("ImportTok")
$importTok.Font.Bold = $true
$importTok.Font.Color = 32768

$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 32768

foreach ($styleName in @("KeywordTok","CommentTok","DocumentationTok","AnnotationTok","CommentVarTok", `
                         "ControlFlowTok","InformationTok","WarningTok","AlertTok","ErrorTok")) {
    $st = $d.Styles($styleName)
    if ($st.Font.Bold -eq $true) { $st.Font.Bold = $true }
    if ($st.Font.Italic -eq $true) { $st.Font.Italic = $true }
}

Write-Output "done"
